$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in hours for "Reading Docs" and "BronchoVision GUI" tasks (August section)
$ws.Range("C40").Value = 1
$ws.Range("C41").Value = 1

# Fill in "Not Paid" total for the August section
$ws.Range("D46").Value = 2

# Update the active selection to reflect where the user left off
$ws.Range("E47").Select()
